$d = $word.ActiveDocument

# Replace "Spiral Model" with "waterfall model" in the body text
# (commit message: Changed 'spiral model' to 'waterfall model.')
$d.Content.Find.Execute("Spiral Model", $true, $true, $false, $false, $false,
                         $true, 1, $false, "waterfall model", 2)
